# Add a new officer registration row (row 3) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "T1234567J"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = "Pending"
$ws.Range("E3").Value = 45769.79456391204
